# Swap the two theme color schemes: the Slide Master's theme (theme1.xml,
# originally the "Integral" / Red Violet colors) becomes the "Office Theme"
# colors, and the Notes Master's theme (theme2.xml, originally "Office
# Theme" colors) becomes the "Integral" / Red Violet colors.

$p = $ppt.ActivePresentation

# Target color values (VBA RGB() BGR-packed integers), in ThemeColorScheme
# item order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
# 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
$integralColors = @(0, 16777215, 5326149, 14473688, 9514467, 13381832, 14460494, 15168839, 14774665, 7555029, 2465643, 9211020)

$masterTcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterTcs.Item($i).RGB = $officeColors[$i - 1]
}

$notesTcs = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesTcs.Item($i).RGB = $integralColors[$i - 1]
}
